$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quiz")

# Update "Marking" (right answer score) total
$ws.Range("B11").Value = 5

# Update "Total" row: correct marks total and the Corr/Total display text
$ws.Range("B12").Value = 110
$ws.Range("E12").Value = "110/140"
